$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(48, 8).Value = 2329.75
$ws.Cells.Item(48, 10).Value = 2329.75
$ws.Cells.Item(48, 12).Value = 6989.25
$ws.Cells.Item(48, 14).Value = -7573.25
$ws.Cells.Item(56, 8).Value = 2329.75
$ws.Cells.Item(56, 10).Value = 2329.75
$ws.Cells.Item(56, 12).Value = 6989.25
$ws.Cells.Item(56, 14).Value = -8057.25
$ws.Cells.Item(111, 8).Value = 929.125
$ws.Cells.Item(111, 9).Value = 975
$ws.Cells.Item(111, 10).Value = 883.25
$ws.Cells.Item(111, 11).Value = 2925
$ws.Cells.Item(111, 12).Value = 2649.75
$ws.Cells.Item(111, 13).Value = 142
$ws.Cells.Item(111, 14).Value = -8783.75
$ws.Cells.Item(112, 8).Value = 8929876
$ws.Cells.Item(112, 10).Value = 9616308
$ws.Cells.Item(112, 12).Value = 28848924
$ws.Cells.Item(112, 14).Value = -28851140
$ws.Cells.Item(113, 8).Value = 5090.909
$ws.Cells.Item(113, 9).Value = 4600
$ws.Cells.Item(113, 11).Value = 4600
$ws.Cells.Item(113, 13).Value = -1346
$ws.Cells.Item(132, 8).Value = 3847703
$ws.Cells.Item(132, 9).Value = 4763192
$ws.Cells.Item(132, 10).Value = 2648.6
$ws.Cells.Item(132, 11).Value = 14289576
$ws.Cells.Item(132, 12).Value = 7945.799999999999
$ws.Cells.Item(132, 13).Value = -14287046
$ws.Cells.Item(132, 14).Value = -13005.8
$ws.Cells.Item(138, 8).Value = 2005.9836
$ws.Cells.Item(138, 9).Value = 1227.921
$ws.Cells.Item(138, 10).Value = 3291.4783
$ws.Cells.Item(138, 11).Value = 3683.763
$ws.Cells.Item(138, 12).Value = 9874.4349
$ws.Cells.Item(138, 13).Value = 1456.237
$ws.Cells.Item(138, 14).Value = -20154.4349
$ws.Cells.Item(141, 8).Value = 677766.9
$ws.Cells.Item(141, 10).Value = 800302.2
$ws.Cells.Item(141, 12).Value = 2400906.6
$ws.Cells.Item(141, 14).Value = -2411266.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6311.891
$ws.Cells.Item(32, 9).Value = 5145.8047
$ws.Cells.Item(32, 10).Value = 26601.8
$ws.Cells.Item(32, 11).Value = 5145.8047
$ws.Cells.Item(32, 12).Value = 26601.8
$ws.Cells.Item(32, 13).Value = -4858.8047
$ws.Cells.Item(32, 14).Value = -27175.8
$ws.Cells.Item(64, 8).Value = 28423.076
$ws.Cells.Item(64, 9).Value = 20000
$ws.Cells.Item(64, 10).Value = 29125
$ws.Cells.Item(64, 11).Value = 20000
$ws.Cells.Item(64, 12).Value = 29125
$ws.Cells.Item(64, 13).Value = -19752
$ws.Cells.Item(64, 14).Value = -29621
$ws.Cells.Item(67, 8).Value = 28423.076
$ws.Cells.Item(67, 9).Value = 20000
$ws.Cells.Item(67, 10).Value = 29125
$ws.Cells.Item(67, 11).Value = 20000
$ws.Cells.Item(67, 12).Value = 29125
$ws.Cells.Item(67, 13).Value = -19142
$ws.Cells.Item(67, 14).Value = -30841
$ws.Cells.Item(97, 8).Value = 545.09375
$ws.Cells.Item(97, 9).Value = 479.32
$ws.Cells.Item(97, 10).Value = 780
$ws.Cells.Item(97, 11).Value = 479.32
$ws.Cells.Item(97, 12).Value = 780
$ws.Cells.Item(97, 13).Value = 16.68000000000001
$ws.Cells.Item(97, 14).Value = -1772
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1467.6957
$ws.Cells.Item(20, 9).Value = 1316.3125
$ws.Cells.Item(20, 10).Value = 1813.7142
$ws.Cells.Item(20, 11).Value = 1316.3125
$ws.Cells.Item(20, 12).Value = 1813.7142
$ws.Cells.Item(20, 13).Value = -1069.3125
$ws.Cells.Item(20, 14).Value = -2307.7142
$ws.Cells.Item(80, 8).Value = 554.5
$ws.Cells.Item(80, 10).Value = 468.25
$ws.Cells.Item(80, 12).Value = 468.25
$ws.Cells.Item(80, 14).Value = -2464.25
$ws.Cells.Item(83, 8).Value = 554.5
$ws.Cells.Item(83, 10).Value = 468.25
$ws.Cells.Item(83, 12).Value = 2341.25
$ws.Cells.Item(83, 14).Value = -12325.25
$ws.Cells.Item(99, 8).Value = 3931.6667
$ws.Cells.Item(99, 9).Value = 1396.6666
$ws.Cells.Item(99, 10).Value = 6466.6665
$ws.Cells.Item(99, 11).Value = 1396.6666
$ws.Cells.Item(99, 12).Value = 6466.6665
$ws.Cells.Item(99, 13).Value = 101.3334
$ws.Cells.Item(99, 14).Value = -9462.666499999999
$ws.Cells.Item(134, 8).Value = 2467.8108
$ws.Cells.Item(134, 9).Value = 1736.2424
$ws.Cells.Item(134, 10).Value = 8503.25
$ws.Cells.Item(134, 11).Value = 5208.7272
$ws.Cells.Item(134, 12).Value = 25509.75
$ws.Cells.Item(134, 13).Value = -2673.7272
$ws.Cells.Item(134, 14).Value = -30579.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1925271
$ws.Cells.Item(31, 9).Value = 2501505.2
$ws.Cells.Item(31, 10).Value = 4489.75
$ws.Cells.Item(31, 11).Value = 2501505.2
$ws.Cells.Item(31, 12).Value = 4489.75
$ws.Cells.Item(31, 13).Value = -2501210.2
$ws.Cells.Item(31, 14).Value = -5079.75
$ws.Cells.Item(34, 8).Value = 1925271
$ws.Cells.Item(34, 9).Value = 2501505.2
$ws.Cells.Item(34, 10).Value = 4489.75
$ws.Cells.Item(34, 11).Value = 2501505.2
$ws.Cells.Item(34, 12).Value = 4489.75
$ws.Cells.Item(34, 13).Value = -2501303.2
$ws.Cells.Item(34, 14).Value = -4893.75
$ws.Cells.Item(52, 8).Value = 25000
$ws.Cells.Item(52, 10).Value = 25000
$ws.Cells.Item(52, 12).Value = 25000
$ws.Cells.Item(52, 14).Value = -25588
$ws.Cells.Item(122, 8).Value = 2009.8372
$ws.Cells.Item(122, 9).Value = 1699.2424
$ws.Cells.Item(122, 10).Value = 3034.8
$ws.Cells.Item(122, 11).Value = 5097.7272
$ws.Cells.Item(122, 12).Value = 9104.400000000001
$ws.Cells.Item(122, 13).Value = -2647.7272
$ws.Cells.Item(122, 14).Value = -14004.4
$ws.Cells.Item(132, 8).Value = 2985.7896
$ws.Cells.Item(132, 9).Value = 2485.8333
$ws.Cells.Item(132, 10).Value = 3842.8572
$ws.Cells.Item(132, 11).Value = 7457.499899999999
$ws.Cells.Item(132, 12).Value = 11528.5716
$ws.Cells.Item(132, 13).Value = -4927.499899999999
$ws.Cells.Item(132, 14).Value = -16588.5716
$ws.Cells.Item(134, 8).Value = 2223.5
$ws.Cells.Item(134, 9).Value = 751.6429000000001
$ws.Cells.Item(134, 11).Value = 2254.9287
$ws.Cells.Item(134, 13).Value = 280.0712999999996
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 58911.39
$ws.Cells.Item(102, 9).Value = 2950.25
$ws.Cells.Item(102, 10).Value = 103680.3
$ws.Cells.Item(102, 11).Value = 2950.25
$ws.Cells.Item(102, 12).Value = 103680.3
$ws.Cells.Item(102, 13).Value = -1328.25
$ws.Cells.Item(102, 14).Value = -106924.3
$ws.Cells.Item(113, 8).Value = 2281
$ws.Cells.Item(113, 9).Value = 1717.2858
$ws.Cells.Item(113, 11).Value = 1717.2858
$ws.Cells.Item(113, 13).Value = 452.7141999999999
$ws.Cells.Item(132, 8).Value = 4077.4
$ws.Cells.Item(132, 9).Value = 4737.4287
$ws.Cells.Item(132, 10).Value = 3499.875
$ws.Cells.Item(132, 11).Value = 14212.2861
$ws.Cells.Item(132, 12).Value = 10499.625
$ws.Cells.Item(132, 13).Value = -11682.2861
$ws.Cells.Item(132, 14).Value = -15559.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4821.722
$ws.Cells.Item(40, 9).Value = 5056.5
$ws.Cells.Item(40, 11).Value = 5056.5
$ws.Cells.Item(40, 13).Value = -4920.5
$ws.Cells.Item(132, 8).Value = 2929
$ws.Cells.Item(132, 9).Value = 2028.4445
$ws.Cells.Item(132, 10).Value = 3405.7646
$ws.Cells.Item(132, 11).Value = 6085.333500000001
$ws.Cells.Item(132, 12).Value = 10217.2938
$ws.Cells.Item(132, 13).Value = -3555.333500000001
$ws.Cells.Item(132, 14).Value = -15277.2938
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 436992.34
$ws.Cells.Item(122, 9).Value = 589871.5
$ws.Cells.Item(122, 10).Value = 3834.8333
$ws.Cells.Item(122, 11).Value = 1769614.5
$ws.Cells.Item(122, 12).Value = 11504.4999
$ws.Cells.Item(122, 13).Value = -1767164.5
$ws.Cells.Item(122, 14).Value = -16404.4999
$ws.Cells.Item(132, 8).Value = 138219.19
$ws.Cells.Item(132, 9).Value = 164573.73
$ws.Cells.Item(132, 10).Value = 4250.25
$ws.Cells.Item(132, 11).Value = 493721.1900000001
$ws.Cells.Item(132, 12).Value = 12750.75
$ws.Cells.Item(132, 13).Value = -491191.1900000001
$ws.Cells.Item(132, 14).Value = -17810.75
